$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the text of G6 (shared string "min error formula" -> "too many formulas")
$ws.Range("G6").Value = "too many formulas"

# Update E6 value from 0.5 to 0.9
$ws.Range("E6").Value = 0.9

# Remove the yellow highlight fill from B6 and C6, keeping their existing number formats
$ws.Range("B6").Interior.Pattern = -4142
$ws.Range("C6").Interior.Pattern = -4142

# Move the active cell selection to D7
$ws.Range("D7").Select()
